$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '85.641.84'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +6.07%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.324.29'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.64%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '219.35'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.61%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '637.25'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.51%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.322'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +10.99%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.595'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.28%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '3.327.56'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.71%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.600'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.44%  '

$ws.Range("E12").Value = '  +1.52%  '

$ws.Range("E13").Value = '  -0.23%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.927.52'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.37%  '

$ws.Range("E15").Value = '  +3.26%  '

$ws.Range("E16").Value = '  -1.50%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '85.126.89'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +5.46%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.313.52'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.80%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '14.64'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '3.19'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.30%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '440.64'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.16%  '

$ws.Range("E22").Value = '  -3.08%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.26'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -2.74%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '7.39'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +4.95%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '5.49'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +12.46%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '12.25'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +10.05%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '3.480.75'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.21%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '78.40'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.48%  '

$ws.Range("E29").Value = '  +1.76%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.05%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '610.08'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +6.59%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.165'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +31.29%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '9.26'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.37%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.03%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.55'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.44%  '

$ws.Range("E36").Value = '  -0.69%  '

$ws.Range("E37").Value = '  -3.11%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '23.27'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.75%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '6.46'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +8.68%  '

$ws.Range("E40").Value = '  -0.74%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '21.23'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +4.28%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.11'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +9.04%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '159.93'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.09%  '

$ws.Range("E46").Value = '  +0.00%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '190.23'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.96%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.37'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.29%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '45.05'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.97%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.791'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.05%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '26.61'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.79%  '
